# "fix: various task transformer improvements" — add new xlsx template
# columns "goal_version", "rule_name_id" and "rule_version" right after the
# existing "goal_name_id" column, re-using the column space (AO:AQ) that was
# already blank between "goal_name_id" (AN) and the "Parameter .../Values
# ..." columns (old AR/AS). Inserting one column ahead of the old AR column
# shifts that pair one slot to the right (-> AS/AT) to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "Parameter .../Values ..." columns one slot to the right.
$ws.Columns("AR").Insert() | Out-Null

# New header cells for row 1 (previously-empty AO1:AQ1).
$ws.Range("AO1").Value = "goal_version"
$ws.Range("AP1").Value = "rule_name_id"
$ws.Range("AQ1").Value = "rule_version"

# Match the formatting of the neighbouring "goal_name_id" header (AN1).
$ws.Range("AN1").Copy() | Out-Null
$ws.Range("AO1:AQ1").PasteSpecial(-4122) | Out-Null

# New (still empty) row-2 cell under "goal_version", styled like AN2.
$ws.Range("AN2").Copy() | Out-Null
$ws.Range("AO2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Leave the selection where the edited workbook left it.
$ws.Range("AN2").Select() | Out-Null
